# Updates exercises, deletes excess files
# Delete the rows for MIND ID 201, 201_FU, 208, 216, 245 from Sheet2's QC table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Delete from the bottom up so row numbers of rows-to-delete above aren't shifted.
$ws.Rows.Item(10).Delete()  # MIND ID 245
$ws.Rows.Item(9).Delete()   # MIND ID 216
$ws.Rows.Item(7).Delete()   # MIND ID 208
$ws.Rows.Item(5).Delete()   # MIND ID 201_FU
$ws.Rows.Item(2).Delete()   # MIND ID 201

$ws.Activate()
$ws.Range("E14").Select()
